# Applies the changes described by the diff:
#  - Removes the fill style (s="1") from B2/B3 on the "Maestro" sheet
#  - Renames the 10 students on the "Desarrollo Web" sheet, replacing
#    the shared strings used by A2:A11, and updates the active cell
#    selection to A12.

$wb = $excel.ActiveWorkbook

# --- Sheet "Maestro": clear the cell style (fill) applied to B2/B3 ---
$maestro = $wb.Worksheets.Item("Maestro")
$maestro.Range("B2:B3").Style = "Normal"

# --- Sheet "Desarrollo Web": update student names ---
$desarrollo = $wb.Worksheets.Item("Desarrollo Web")

$newNames = @(
    "Gael Barroso",
    "Maria Ines Vargas",
    "Enriqueta Pérez",
    "Carmelo Gascon",
    "Ian Lago",
    "Fatiha Agullo",
    "Cristina Maria Prados",
    "Jairo Vila",
    "Maria Alicia Roca",
    "Hipolito Montiel"
)

for ($i = 0; $i -lt $newNames.Length; $i++) {
    $row = $i + 2
    $desarrollo.Cells.Item($row, 1).Value = $newNames[$i]
}

$desarrollo.Activate()
$desarrollo.Range("A12").Select()
